$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Il1a"
$ws.Cells.Item(2, 3).Value = "Il1rap"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.02286966666666667
$ws.Cells.Item(2, 8).Value = 0.068609
$ws.Cells.Item(2, 9).Value = 0.001711767187487096
$ws.Cells.Item(2, 10).Value = 0.001711767187487096
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 3.240187666666666
$ws.Cells.Item(2, 14).Value = 9.720562999999999
$ws.Cells.Item(2, 15).Value = 0.09865275843473079
$ws.Cells.Item(2, 16).Value = 0.1091454159637843
$ws.Cells.Item(2, 17).Value = 0.0741020118741111
$ws.Cells.Item(2, 18).Value = 0.6669181068669999
$ws.Cells.Item(2, 19).Value = 0.000168870554843663
$ws.Cells.Item(2, 20).Value = 0.0001868315417114363
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Il1a"
$ws.Cells.Item(3, 3).Value = "Il1rap"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.02286966666666667
$ws.Cells.Item(3, 8).Value = 0.068609
$ws.Cells.Item(3, 9).Value = 0.001711767187487096
$ws.Cells.Item(3, 10).Value = 0.001711767187487096
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 8.637706333333334
$ws.Cells.Item(3, 14).Value = 25.913119
$ws.Cells.Item(3, 15).Value = 0.2629889512569831
$ws.Cells.Item(3, 16).Value = 0.2909603232008314
$ws.Cells.Item(3, 17).Value = 0.1975414646078889
$ws.Cells.Item(3, 18).Value = 1.777873181471
$ws.Cells.Item(3, 19).Value = 0.0004501758574333471
$ws.Cells.Item(3, 20).Value = 0.0004980563341158237
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Il1a"
$ws.Cells.Item(4, 3).Value = "Il1rap"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.02286966666666667
$ws.Cells.Item(4, 8).Value = 0.068609
$ws.Cells.Item(4, 9).Value = 0.001711767187487096
$ws.Cells.Item(4, 10).Value = 0.001711767187487096
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 3.343997
$ws.Cells.Item(4, 14).Value = 10.031991
$ws.Cells.Item(4, 15).Value = 0.1018134016252344
$ws.Cells.Item(4, 16).Value = 0.1126422235666741
$ws.Cells.Item(4, 17).Value = 0.07647609672433334
$ws.Cells.Item(4, 18).Value = 0.688284870519
$ws.Cells.Item(4, 19).Value = 0.0001742808401485217
$ws.Cells.Item(4, 20).Value = 0.0001928172622270185
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Il1a"
$ws.Cells.Item(5, 3).Value = "Il1rap"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.02286966666666667
$ws.Cells.Item(5, 8).Value = 0.068609
$ws.Cells.Item(5, 9).Value = 0.001711767187487096
$ws.Cells.Item(5, 10).Value = 0.001711767187487096
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.150031666666665
$ws.Cells.Item(5, 14).Value = 24.450095
$ws.Cells.Item(5, 15).Value = 0.2481409066266243
$ws.Cells.Item(5, 16).Value = 0.2745330480476329
$ws.Cells.Item(5, 17).Value = 0.1863885075394444
$ws.Cells.Item(5, 18).Value = 1.677496567855
$ws.Cells.Item(5, 19).Value = 0.0004247594618367549
$ws.Cells.Item(5, 20).Value = 0.0004699366635287565
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Il1a"
$ws.Cells.Item(6, 3).Value = "Il1rap"
$ws.Cells.Item(6, 4).Value = "sCs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.02286966666666667
$ws.Cells.Item(6, 8).Value = 0.068609
$ws.Cells.Item(6, 9).Value = 0.001711767187487096
$ws.Cells.Item(6, 10).Value = 0.001711767187487096
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 9.472446999999999
$ws.Cells.Item(6, 14).Value = 18.944894
$ws.Cells.Item(6, 15).Value = 0.2884039820564273
$ws.Cells.Item(6, 16).Value = 0.2127189892210772
$ws.Cells.Item(6, 17).Value = 0.2166317054076666
$ws.Cells.Item(6, 18).Value = 1.299790232446
$ws.Cells.Item(6, 19).Value = 0.0004936804732248095
$ws.Cells.Item(6, 20).Value = 0.0003641253859040612
$ws.Cells.Item(7, 1).Value = "M1"
$ws.Cells.Item(7, 2).Value = "Il1a"
$ws.Cells.Item(7, 3).Value = "Il1rap"
$ws.Cells.Item(7, 4).Value = "ECs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 5.046308333333333
$ws.Cells.Item(7, 8).Value = 15.138925
$ws.Cells.Item(7, 9).Value = 0.3777101410722805
$ws.Cells.Item(7, 10).Value = 0.3777101410722805
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 3.240187666666666
$ws.Cells.Item(7, 14).Value = 9.720562999999999
$ws.Cells.Item(7, 15).Value = 0.09865275843473079
$ws.Cells.Item(7, 16).Value = 0.1091454159637843
$ws.Cells.Item(7, 17).Value = 16.35098602386389
$ws.Cells.Item(7, 18).Value = 147.158874214775
$ws.Cells.Item(7, 19).Value = 0.03726214730555177
$ws.Cells.Item(7, 20).Value = 0.0412253304610737
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 2).Value = "Il1a"
$ws.Cells.Item(8, 3).Value = "Il1rap"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.046308333333333
$ws.Cells.Item(8, 8).Value = 15.138925
$ws.Cells.Item(8, 9).Value = 0.3777101410722805
$ws.Cells.Item(8, 10).Value = 0.3777101410722805
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.637706333333334
$ws.Cells.Item(8, 14).Value = 25.913119
$ws.Cells.Item(8, 15).Value = 0.2629889512569831
$ws.Cells.Item(8, 16).Value = 0.2909603232008314
$ws.Cells.Item(8, 17).Value = 43.58852945078611
$ws.Cells.Item(8, 18).Value = 392.296765057075
$ws.Cells.Item(8, 19).Value = 0.09933359387972619
$ws.Cells.Item(8, 20).Value = 0.1098986647226223
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 2).Value = "Il1a"
$ws.Cells.Item(9, 3).Value = "Il1rap"
$ws.Cells.Item(9, 4).Value = "M1"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.046308333333333
$ws.Cells.Item(9, 8).Value = 15.138925
$ws.Cells.Item(9, 9).Value = 0.3777101410722805
$ws.Cells.Item(9, 10).Value = 0.3777101410722805
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.343997
$ws.Cells.Item(9, 14).Value = 10.031991
$ws.Cells.Item(9, 15).Value = 0.1018134016252344
$ws.Cells.Item(9, 16).Value = 0.1126422235666741
$ws.Cells.Item(9, 17).Value = 16.87483992774166
$ws.Cells.Item(9, 18).Value = 151.873559349675
$ws.Cells.Item(9, 19).Value = 0.03845595429091604
$ws.Cells.Item(9, 20).Value = 0.04254611015406384
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 2).Value = "Il1a"
$ws.Cells.Item(10, 3).Value = "Il1rap"
$ws.Cells.Item(10, 4).Value = "M2"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.046308333333333
$ws.Cells.Item(10, 8).Value = 15.138925
$ws.Cells.Item(10, 9).Value = 0.3777101410722805
$ws.Cells.Item(10, 10).Value = 0.3777101410722805
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 8.150031666666665
$ws.Cells.Item(10, 14).Value = 24.450095
$ws.Cells.Item(10, 15).Value = 0.2481409066266243
$ws.Cells.Item(10, 16).Value = 0.2745330480476329
$ws.Cells.Item(10, 17).Value = 41.12757271643055
$ws.Cells.Item(10, 18).Value = 370.148154447875
$ws.Cells.Item(10, 19).Value = 0.09372533684774585
$ws.Cells.Item(10, 20).Value = 0.1036939163070746
$ws.Cells.Item(11, 1).Value = "M1"
$ws.Cells.Item(11, 2).Value = "Il1a"
$ws.Cells.Item(11, 3).Value = "Il1rap"
$ws.Cells.Item(11, 4).Value = "sCs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 5.046308333333333
$ws.Cells.Item(11, 8).Value = 15.138925
$ws.Cells.Item(11, 9).Value = 0.3777101410722805
$ws.Cells.Item(11, 10).Value = 0.3777101410722805
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 9.472446999999999
$ws.Cells.Item(11, 14).Value = 18.944894
$ws.Cells.Item(11, 15).Value = 0.2884039820564273
$ws.Cells.Item(11, 16).Value = 0.2127189892210772
$ws.Cells.Item(11, 17).Value = 47.80088823315833
$ws.Cells.Item(11, 18).Value = 286.80532939895
$ws.Cells.Item(11, 19).Value = 0.1089331087483406
$ws.Cells.Item(11, 20).Value = 0.08034611942744596
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Il1a"
$ws.Cells.Item(12, 3).Value = "Il1rap"
$ws.Cells.Item(12, 4).Value = "ECs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 8.291089
$ws.Cells.Item(12, 8).Value = 24.873267
$ws.Cells.Item(12, 9).Value = 0.6205780917402324
$ws.Cells.Item(12, 10).Value = 0.6205780917402324
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.240187666666666
$ws.Cells.Item(12, 14).Value = 9.720562999999999
$ws.Cells.Item(12, 15).Value = 0.09865275843473079
$ws.Cells.Item(12, 16).Value = 0.1091454159637843
$ws.Cells.Item(12, 17).Value = 26.86468432103566
$ws.Cells.Item(12, 18).Value = 241.782158889321
$ws.Cells.Item(12, 19).Value = 0.06122174057433535
$ws.Cells.Item(12, 20).Value = 0.06773325396099915
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Il1a"
$ws.Cells.Item(13, 3).Value = "Il1rap"
$ws.Cells.Item(13, 4).Value = "FAPs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 8.291089
$ws.Cells.Item(13, 8).Value = 24.873267
$ws.Cells.Item(13, 9).Value = 0.6205780917402324
$ws.Cells.Item(13, 10).Value = 0.6205780917402324
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 8.637706333333334
$ws.Cells.Item(13, 14).Value = 25.913119
$ws.Cells.Item(13, 15).Value = 0.2629889512569831
$ws.Cells.Item(13, 16).Value = 0.2909603232008314
$ws.Cells.Item(13, 17).Value = 71.61599196553033
$ws.Cells.Item(13, 18).Value = 644.543927689773
$ws.Cells.Item(13, 19).Value = 0.1632051815198236
$ws.Cells.Item(13, 20).Value = 0.1805636021440932
$ws.Cells.Item(14, 1).Value = "M2"
$ws.Cells.Item(14, 2).Value = "Il1a"
$ws.Cells.Item(14, 3).Value = "Il1rap"
$ws.Cells.Item(14, 4).Value = "M1"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 8.291089
$ws.Cells.Item(14, 8).Value = 24.873267
$ws.Cells.Item(14, 9).Value = 0.6205780917402324
$ws.Cells.Item(14, 10).Value = 0.6205780917402324
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 3.343997
$ws.Cells.Item(14, 14).Value = 10.031991
$ws.Cells.Item(14, 15).Value = 0.1018134016252344
$ws.Cells.Item(14, 16).Value = 0.1126422235666741
$ws.Cells.Item(14, 17).Value = 27.725376742733
$ws.Cells.Item(14, 18).Value = 249.528390684597
$ws.Cells.Item(14, 19).Value = 0.06318316649416984
$ws.Cells.Item(14, 20).Value = 0.06990329615038327
$ws.Cells.Item(15, 1).Value = "M2"
$ws.Cells.Item(15, 2).Value = "Il1a"
$ws.Cells.Item(15, 3).Value = "Il1rap"
$ws.Cells.Item(15, 4).Value = "M2"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 8.291089
$ws.Cells.Item(15, 8).Value = 24.873267
$ws.Cells.Item(15, 9).Value = 0.6205780917402324
$ws.Cells.Item(15, 10).Value = 0.6205780917402324
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 8.150031666666665
$ws.Cells.Item(15, 14).Value = 24.450095
$ws.Cells.Item(15, 15).Value = 0.2481409066266243
$ws.Cells.Item(15, 16).Value = 0.2745330480476329
$ws.Cells.Item(15, 17).Value = 67.57263790115165
$ws.Cells.Item(15, 18).Value = 608.1537411103649
$ws.Cells.Item(15, 19).Value = 0.1539908103170417
$ws.Cells.Item(15, 20).Value = 0.1703691950770296
$ws.Cells.Item(16, 1).Value = "M2"
$ws.Cells.Item(16, 2).Value = "Il1a"
$ws.Cells.Item(16, 3).Value = "Il1rap"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 8.291089
$ws.Cells.Item(16, 8).Value = 24.873267
$ws.Cells.Item(16, 9).Value = 0.6205780917402324
$ws.Cells.Item(16, 10).Value = 0.6205780917402324
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 9.472446999999999
$ws.Cells.Item(16, 14).Value = 18.944894
$ws.Cells.Item(16, 15).Value = 0.2884039820564273
$ws.Cells.Item(16, 16).Value = 0.2127189892210772
$ws.Cells.Item(16, 17).Value = 78.53690112478299
$ws.Cells.Item(16, 18).Value = 471.2214067486979
$ws.Cells.Item(16, 19).Value = 0.1789771928348619
$ws.Cells.Item(16, 20).Value = 0.1320087444077271
